$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "Email"
$ws.Range("H2").Value = "adrianrentea01@gmail.com"

$ws.Range("H1:H2").Select()
$ws.Columns.Item(8).ColumnWidth = 23.666666666666668
